$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 91
$ws.Range("G6").Value = 2719.08
$ws.Range("B10").Value = 29487.59
$ws.Range("F31").Value = 30
$ws.Range("G31").Value = 1059.6
$ws.Range("B32").Value = 13009.49
$ws.Range("F70").Value = 23
$ws.Range("G70").Value = 3103.85
$ws.Range("F71").Value = 346
$ws.Range("G71").Value = 22040.2
$ws.Range("B90").Value = 192304.64
$ws.Range("F144").Value = 1148
$ws.Range("G144").Value = 9700.6
$ws.Range("B147").Value = 17049.15
$ws.Range("F150").Value = 42
$ws.Range("G150").Value = 1952.58
$ws.Range("B156").Value = 33948.03
$ws.Range("F164").Value = 68
$ws.Range("G164").Value = 7794.84
$ws.Range("B175").Value = 32662.19
$ws.Range("F187").Value = 23
$ws.Range("G187").Value = 1149.31
$ws.Range("F190").Value = 6
$ws.Range("G190").Value = 492.06
$ws.Range("B192").Value = 48706
$ws.Range("E192").Value = 39.8
$ws.Range("F192").Value = -144
$ws.Range("G192").Value = -4795.2
$ws.Range("B193").Value = 64973
$ws.Range("E193").Value = 35.4
$ws.Range("F193").Value = 2
$ws.Range("G193").Value = 66.59999999999999
$ws.Range("F203").Value = 64
$ws.Range("G203").Value = 1290.24
$ws.Range("B216").Value = 45568.95
$ws.Range("F225").Value = 77
$ws.Range("G225").Value = 8795.709999999999
$ws.Range("B227").Value = 55373
$ws.Range("E227").Value = 163.62
$ws.Range("F227").Value = -94
$ws.Range("G227").Value = -13562.32
$ws.Range("B228").Value = 63520
$ws.Range("E228").Value = 153.4
$ws.Range("F228").Value = 67
$ws.Range("G228").Value = 9666.76
$ws.Range("F255").Value = 593
$ws.Range("G255").Value = 101598.69
$ws.Range("B260").Value = 200584.43
$ws.Range("F277").Value = 2
$ws.Range("G277").Value = 42.5
$ws.Range("F278").Value = 14
$ws.Range("G278").Value = 1922.48
$ws.Range("F283").Value = 41
$ws.Range("G283").Value = 14000.27
$ws.Range("F285").Value = 5
$ws.Range("G285").Value = 139.65
$ws.Range("F295").Value = 3
$ws.Range("G295").Value = 311.07
$ws.Range("F302").Value = 67
$ws.Range("G302").Value = 14129.63
$ws.Range("B304").Value = 186563.83
$ws.Range("F321").Value = 44
$ws.Range("G321").Value = 2416.48
$ws.Range("B322").Value = 58047
$ws.Range("D322").Value = 105.54
$ws.Range("E322").Value = 126.1
$ws.Range("F322").Value = 40
$ws.Range("G322").Value = 4221.6
$ws.Range("B323").Value = 47097
$ws.Range("D323").Value = 112.28
$ws.Range("E323").Value = 134.16
$ws.Range("F323").Value = 15
$ws.Range("G323").Value = 1684.2
$ws.Range("F328").Value = 48
$ws.Range("G328").Value = 1786.08
$ws.Range("B330").Value = 29720.33
$ws.Range("F355").Value = 14
$ws.Range("G355").Value = 2288.3
$ws.Range("F357").Value = 6
$ws.Range("G357").Value = 1567.8
$ws.Range("B358").Value = 36761.26
$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9
$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29
$ws.Range("B375").Value = 45718
$ws.Range("E375").Value = 19.38
$ws.Range("F375").Value = -294
$ws.Range("G375").Value = -4768.68
$ws.Range("B376").Value = 64927
$ws.Range("E376").Value = 17.26
$ws.Range("F376").Value = 106
$ws.Range("G376").Value = 1719.32
$ws.Range("F406").Value = 3
$ws.Range("G406").Value = 22.68
$ws.Range("F410").Value = 0
$ws.Range("G410").Value = 0
$ws.Range("B411").Value = 7851.26
$ws.Range("F432").Value = 0
$ws.Range("G432").Value = 0
$ws.Range("F434").Value = 26
$ws.Range("G434").Value = 848.64
$ws.Range("B435").Value = 1019.9
$ws.Range("B442").Value = 64810
$ws.Range("E442").Value = 291.22
$ws.Range("F442").Value = 4
$ws.Range("G442").Value = 1095.68
$ws.Range("B443").Value = 53319
$ws.Range("E443").Value = 310.64
$ws.Range("F443").Value = -6
$ws.Range("G443").Value = -1643.52
$ws.Range("F450").Value = 13
$ws.Range("G450").Value = 1803.62
$ws.Range("F452").Value = 4
$ws.Range("G452").Value = 766.72
$ws.Range("B460").Value = 14306.67
$ws.Range("F485").Value = 23
$ws.Range("G485").Value = 4035.81
$ws.Range("B488").Value = 31973.81
$ws.Range("F509").Value = 237
$ws.Range("G509").Value = 19050.06
$ws.Range("B510").Value = 25182.52
$ws.Range("F550").Value = 4
$ws.Range("G550").Value = 326.24
$ws.Range("F552").Value = 21
$ws.Range("G552").Value = 2137.59
$ws.Range("B560").Value = 6873.69
$ws.Range("B572").Value = 65362
$ws.Range("F572").Value = 24
$ws.Range("G572").Value = 980.88
$ws.Range("B573").Value = 65079
$ws.Range("F573").Value = 18
$ws.Range("G573").Value = 735.66
$ws.Range("F579").Value = 36
$ws.Range("G579").Value = 2901.6
$ws.Range("F580").Value = 64
$ws.Range("G580").Value = 3647.36
$ws.Range("F581").Value = 8
$ws.Range("G581").Value = 1934.4
$ws.Range("B583").Value = 23619.85
$ws.Range("F599").Value = 1872
$ws.Range("G599").Value = 305341.92
$ws.Range("F601").Value = 445
$ws.Range("G601").Value = 125877.15
$ws.Range("B606").Value = 481537.42
$ws.Range("F613").Value = 143
$ws.Range("G613").Value = 22759.88
$ws.Range("B618").Value = 45203.3
$ws.Range("B619").Value = 1885645.24
$ws.Range("B620").Value = 1885645.24
